# Changes of DEV URL configuration
# Update tracking-number / rate / result values on Sheet1.
# These values look numeric (or currency) but must be stored as TEXT
# (shared-string) cells, matching the original column's cell type.
# Trick: force the "@" text number format before assigning the value so
# Excel doesn't auto-coerce it to a number, then reapply the "Normal"
# cell style so no stray number-format / border override is left behind
# on the cell (matches the un-styled <c t="s"> cells in the target file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2: P2 tracking number changes
Set-TextValue "P2" "320018787850"

# Row 3: P3 tracking number changes
Set-TextValue "P3" "320018787860"

# Row 4: P4 tracking number changes
Set-TextValue "P4" "320018764881"

# Row 5: P5/Q5/R5 were empty (bordered) cells, now populated with a
# tracking number, an actual rate, and a PASS/FAIL result.
Set-TextValue "P5" "320018766656"
Set-TextValue "Q5" "`$43.36"
Set-TextValue "R5" "FAIL"
